$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A slightly (13.85546875 -> 14.85546875 raw units)
$ws.Columns.Item(1).ColumnWidth = 14

# New data rows appended below the existing data (row 2 is the last used row)
$rows = @(
    @{ Row = 3; Date = 42600.786759259259; Method = "Bag"; C = 1525; D = 3071; E = 356; F = 62;  G = 13; H = 82; I = 17; J = 0; K = 0; L = 0; M = 0 },
    @{ Row = 4; Date = 42600.825416666667; Method = "Bag"; C = 1296; D = 2515; E = 318; F = 31;  G = 6;  H = 83; I = 16; J = 0; K = 0; L = 0; M = 0 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Date
    $ws.Cells.Item($row, 2).Value = $r.Method
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
}

# Reuse the existing date-formatted style (column A, row 2) for the new A cells
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Range("A3:A4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
